# Insert a brand-new weekly data row for Orégano before the current row 91,
# pushing every following record (old rows 91-143) down by one row
# (new rows 92-144). Excel's Rows(...).Insert() shifts values, formulas and
# formatting automatically, including extending the used range /
# dimension from A1:R143 to A1:R144 and carrying the date number format
# from row 90 into the freshly inserted row 91 (column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(91).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(91, 1).Value = 6
$ws.Cells.Item(91, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(91, 3).Value = "Metropolitana"
$ws.Cells.Item(91, 4).Value = 44574
$ws.Cells.Item(91, 5).Value = 13
$ws.Cells.Item(91, 6).Value = 100112029
$ws.Cells.Item(91, 7).Value = "Orégano"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 34
$ws.Cells.Item(91, 11).Value = 8000
$ws.Cells.Item(91, 12).Value = 9000
$ws.Cells.Item(91, 13).Value = 8441
$ws.Cells.Item(91, 14).Value = "`$/docena de atados"
$ws.Cells.Item(91, 15).Value = "Región Metropolitana"
$ws.Cells.Item(91, 16).Value = 2814
$ws.Cells.Item(91, 17).Value = 3
$ws.Cells.Item(91, 18).Value = "Hortaliza"
